$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (field_model) values for the data rows all become "Image"
$ws.Range("D2:D7").Value = "Image"

# Update the active selection to D7
$ws.Range("D7").Select()
